# Add data for 2022-06-29 (extend "through June 20" to "through June 21")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the header label for the rolling "June 2022" column
$ws.Name = "Through 2022-06-21"
$ws.Range("B1").Value = "June 2022 (through June 21)"

# Helper to add to an existing value (or set it if currently blank)
function Add-Value($cell, [double]$amount) {
    $cur = $ws.Range($cell).Value()
    if ($cur -eq $null -or $cur -eq "") {
        $ws.Range($cell).Value = $amount
    } else {
        $ws.Range($cell).Value = [double]$cur + $amount
    }
}

# North Lawndale (row 4)
Add-Value "AF4" 1

# South Shore (row 5)
Add-Value "AF5" 1

# Humboldt Park (row 6)
Add-Value "B6" 1

# Grand Boulevard (row 7)
Add-Value "H7" 1
Add-Value "N7" 1

# Garfield Park (row 10)
Add-Value "B10" 1
Add-Value "N10" 1
Add-Value "AF10" 1

# Austin (row 14)
Add-Value "H14" 1

# Chicago Lawn (row 24)
Add-Value "Z24" 1

# Ukrainian Village (row 28)
Add-Value "Z28" 1

# Morgan Park (row 35)
Add-Value "H35" 1

# Portage Park (row 82)
Add-Value "H82" 1
